# Insert a new weekly record at row 91 (pushing the existing rows 91-115
# down to 92-116), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row above the current row 91; Excel shifts rows 91:115
# down to 92:116 and the used range grows to A1:R116.
$ws.Rows("91:91").Insert(-4121)

# Populate the freshly inserted row 91 with the new weekly reading.
$ws.Range("A91").Value = 5
$ws.Range("B91").Value = "Macroferia Regional de Talca"
$ws.Range("C91").Value = "Maule"
$ws.Range("D91").Value = 45215
$ws.Range("E91").Value = 7
$ws.Range("F91").Value = 300000000
$ws.Range("G91").Value = "Espárragos"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 1200
$ws.Range("M91").Value = 1200
$ws.Range("N91").Value = "$/kilo"
$ws.Range("O91").Value = "Región del Maule"
$ws.Range("P91").Value = 1200
$ws.Range("Q91").Value = 1
$ws.Range("R91").Value = "Hortaliza"
